$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.394.60"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "3.483.73"
$ws.Range("E3").Value = "  -2.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "188.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.210"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.646"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.33%  "
$ws.Range("E12").Value = "  -3.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").Value = "4.040.90"
$ws.Range("E14").Value = "  -2.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "610.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.82%  "
$ws.Range("D16").Value = "69.478.69"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.475.46"
$ws.Range("E19").Value = "  -2.04%  "
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.981"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "107.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.64%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.75%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.71%  "
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = "  -6.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").Value = "3.607.70"
$ws.Range("E37").Value = "  -0.47%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.38%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.392"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "507.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.65%  "
$ws.Range("D42").Value = "0.0₃0767"
$ws.Range("E42").Value = "  -6.88%  "
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("E44").Value = "  -3.16%  "
$ws.Range("E45").Value = "  -2.37%  "
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("E47").Value = "  -4.44%  "
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.16%  "
$ws.Range("E51").Value = "  -7.62%  "
